# Implementado el WhatsApp web
# Adds "WhatsApp Profesional" / "WhatsApp Tramitador" columns to both
# worksheets, plus three new rows of data on "Obras en general", and
# clears several placeholder-empty cells that previously held the
# "inlineStr"-typed blank values in row 5.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Obras en general"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# New bold header cells in columns Y/Z
$ws1.Range("Y1").Value = "WhatsApp Profesional"
$ws1.Range("Y1").Font.Bold = $true
$ws1.Range("Z1").Value = "WhatsApp Tramitador"
$ws1.Range("Z1").Font.Bold = $true

# Row 5 had several cells that used to be touched-but-blank; clear them
# out completely so they disappear from the sheet, matching the rest of
# the already-blank columns in that row.
$ws1.Range("D5").ClearContents()
$ws1.Range("L5:Q5").ClearContents()
$ws1.Range("T5:W5").ClearContents()

# New row 6
$ws1.Range("A6").Value = "'02/06/2025"
$ws1.Range("A6").Style = "Normal"
$ws1.Range("B6").Value = "Ingeniero"
$ws1.Range("C6").Value = "Físico"
$ws1.Range("D6").Value = "'3"
$ws1.Range("D6").Style = "Normal"
$ws1.Range("E6").Value = "Obra nueva"
$ws1.Range("F6").Value = "GUSTI AGUINAGALDE"
$ws1.Range("G6").Value = "AGUSTIN AGUINAGALDE"
$ws1.Range("H6").Value = "FRANCISCO DE HARO 27485"
$ws1.Range("I6").Value = "1516/H/515"
$ws1.Range("K6").Value = "'151515"
$ws1.Range("K6").Style = "Normal"
$ws1.Range("L6").Value = "'7000"
$ws1.Range("L6").Style = "Normal"
$ws1.Range("N6").Value = "'5000"
$ws1.Range("N6").Style = "Normal"
$ws1.Range("O6").Value = "'9000"
$ws1.Range("O6").Style = "Normal"
$ws1.Range("R6").Value = "No pagado"
$ws1.Range("S6").Value = "No pagado"
$ws1.Range("Y6").Value = "'3764251817"
$ws1.Range("Y6").Style = "Normal"
$ws1.Range("Z6").Value = "'3764251817"
$ws1.Range("Z6").Style = "Normal"

# New row 7
$ws1.Range("A7").Value = "'02/06/2025"
$ws1.Range("A7").Style = "Normal"
$ws1.Range("B7").Value = "Ingeniero"
$ws1.Range("C7").Value = "Físico"
$ws1.Range("D7").Value = "'4"
$ws1.Range("D7").Style = "Normal"
$ws1.Range("E7").Value = "Registración"
$ws1.Range("F7").Value = "mario alberto"
$ws1.Range("G7").Value = "lucia benitez"
$ws1.Range("H7").Value = "quesreyooo"
$ws1.Range("I7").Value = "151/h/18198"
$ws1.Range("K7").Value = "'891891"
$ws1.Range("K7").Style = "Normal"
$ws1.Range("L7").Value = "'15000"
$ws1.Range("L7").Style = "Normal"
$ws1.Range("N7").Value = "'999"
$ws1.Range("N7").Style = "Normal"
$ws1.Range("O7").Value = "'999"
$ws1.Range("O7").Style = "Normal"
$ws1.Range("P7").Value = "'100"
$ws1.Range("P7").Style = "Normal"
$ws1.Range("R7").Value = "No pagado"
$ws1.Range("S7").Value = "No pagado"
$ws1.Range("Y7").Value = "'3764251817"
$ws1.Range("Y7").Style = "Normal"
$ws1.Range("Z7").Value = "'3765100571"
$ws1.Range("Z7").Style = "Normal"

# New row 8
$ws1.Range("A8").Value = "'02/06/2025"
$ws1.Range("A8").Style = "Normal"
$ws1.Range("B8").Value = "MMO"
$ws1.Range("C8").Value = "Físico"
$ws1.Range("D8").Value = "'5"
$ws1.Range("D8").Style = "Normal"
$ws1.Range("E8").Value = "Registración"
$ws1.Range("F8").Value = "RIOS DANIEL"
$ws1.Range("G8").Value = "GUSTAVITO"
$ws1.Range("H8").Value = "FRANCISOC DE HARO 2745"
$ws1.Range("I8").Value = "151/1851"
$ws1.Range("K8").Value = "'15165"
$ws1.Range("K8").Style = "Normal"
$ws1.Range("L8").Value = "'70000"
$ws1.Range("L8").Style = "Normal"
$ws1.Range("N8").Value = "'10000"
$ws1.Range("N8").Style = "Normal"
$ws1.Range("O8").Value = "'12000"
$ws1.Range("O8").Style = "Normal"
$ws1.Range("P8").Value = "'10000"
$ws1.Range("P8").Style = "Normal"
$ws1.Range("Q8").Value = "'"
$ws1.Range("Q8").Style = "Normal"
$ws1.Range("R8").Value = "No pagado"
$ws1.Range("S8").Value = "No pagado"
$ws1.Range("T8").Value = "'"
$ws1.Range("T8").Style = "Normal"
$ws1.Range("U8").Value = "'"
$ws1.Range("U8").Style = "Normal"
$ws1.Range("V8").Value = "'"
$ws1.Range("V8").Style = "Normal"
$ws1.Range("W8").Value = "'"
$ws1.Range("W8").Style = "Normal"
$ws1.Range("Y8").Value = "'3765100571"
$ws1.Range("Y8").Style = "Normal"
$ws1.Range("Z8").Value = "'3764251817"
$ws1.Range("Z8").Style = "Normal"

# ---------------------------------------------------------------------
# Sheet 2: "Informes técnicos"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("P1").Value = "WhatsApp Profesional"
$ws2.Range("P1").Font.Bold = $true
$ws2.Range("Q1").Value = "WhatsApp Tramitador"
$ws2.Range("Q1").Font.Bold = $true
